$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "1W"
$ws.Cells.Item(2, 2).Value = 0.024
$ws.Cells.Item(3, 1).Value = "4C"
$ws.Cells.Item(3, 2).Value = 0.994
$ws.Cells.Item(4, 1).Value = "4I"
$ws.Cells.Item(4, 2).Value = 1.015666667
$ws.Cells.Item(5, 1).Value = "4W"
$ws.Cells.Item(5, 2).Value = 0.666
$ws.Cells.Item(6, 1).Value = "5C"
$ws.Cells.Item(6, 2).Value = 1.0095
$ws.Cells.Item(7, 1).Value = "5F"
$ws.Cells.Item(7, 2).Value = 1.0485
$ws.Cells.Item(8, 1).Value = "5I"
$ws.Cells.Item(8, 2).Value = 1.049
$ws.Cells.Item(9, 1).Value = "5L"
$ws.Cells.Item(9, 2).Value = 1.0235
$ws.Cells.Item(10, 1).Value = "5N"
$ws.Cells.Item(10, 2).Value = 1.0465
$ws.Cells.Item(11, 1).Value = "5V"
$ws.Cells.Item(11, 2).Value = 1.0095
$ws.Cells.Item(12, 1).Value = "5Y"
$ws.Cells.Item(12, 2).Value = 1.0605
$ws.Cells.Item(13, 1).Value = "25Q"
$ws.Cells.Item(13, 2).Value = 0.994
